$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price / 1h-volume figures (and, for rows 48-50, the coin name/link
# as the ranking shuffled) to reflect the latest scrape.
# For D-column values that look like plain numbers, force a text number
# format first so Excel keeps the exact original text (e.g. trailing
# zeros, thousands separators) instead of silently converting to a number.
$ws.Range("D2").Value = "43.986.62"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.353.85"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +5.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.71"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.50"
$ws.Range("E7").Value = "  +3.72%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +20.90%  "
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.36"
$ws.Range("E12").Value = "  +22.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("E13").Value = "  +13.34%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "2.704.53"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.81"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.930"
$ws.Range("E17").Value = "  +5.83%  "
$ws.Range("D18").Value = "2.351.85"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "43.822.89"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.70"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "262.96"
$ws.Range("E23").Value = "  +5.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("E26").Value = "  -5.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  +7.31%  "
$ws.Range("E28").Value = "  +17.62%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.12"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.38"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("E35").Value = "  +9.39%  "
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.42"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0282"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.220"
$ws.Range("E41").Value = "  +22.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.22"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.15"
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.106"
$ws.Range("E44").Value = "  +10.08%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  +11.07%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.62"
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.89"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.19"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.30"
$ws.Range("E51").Value = "  +10.56%  "

# Restore default style for cells where we forced text number format,
# so they match the original (unstyled) appearance.
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
